# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# (and the blank paragraph that separated it from the bibliography entry
# above it) that the Jekyll site rebuild dropped from the page.

$d = $word.ActiveDocument

# Locate the last bibliography paragraph ("2000.SCHREIBER, G,.P. - Usinas
# Hidrelétricas ...") that the footer block immediately follows.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2000.SCHREIBER*") {
        $anchor = $p
    }
}

if ($anchor -ne $null) {
    # The next three paragraphs are: an empty spacer paragraph, the
    # "Ver no Jupiter Salvar em pdf Salvar em docx" line, and the
    # "© 2020 . Contact: ..." credit line. Delete all three as one
    # range so the remaining empty/page-break paragraphs after them
    # are left untouched.
    $first = $anchor.Next()
    $second = $first.Next()
    $third = $second.Next()

    $deadRange = $d.Range($first.Range.Start, $third.Range.End)
    $deadRange.Delete()
}
